$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.040.64"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.789.90"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.79%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.01"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.07%  "

$ws.Range("E6").Value = "  +0.93%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.10"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.54%  "

$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("E10").Value = "  -1.37%  "

$ws.Range("E11").Value = "  -3.38%  "

$ws.Range("E12").Value = "  -0.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.048.67"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.31"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +9.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.779.30"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.38%  "

$ws.Range("E16").Value = "  -1.83%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.055.99"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.22"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "253.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.13%  "

$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("E24").Value = "  -2.93%  "

$ws.Range("E25").Value = "  -2.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.32"
$ws.Range("D26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.86%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("E29").Value = "  -2.75%  "

$ws.Range("E30").Value = "  +0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.15%  "

$ws.Range("E32").Value = "  -0.09%  "

$ws.Range("E33").Value = "  -1.82%  "

$ws.Range("E34").Value = "  +1.42%  "

$ws.Range("E35").Value = "  -0.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.467.89"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.64%  "

$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.633"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("E39").Value = "  -1.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.43%  "

$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.902"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("E44").Value = "  -2.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0514"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.67%  "

$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.947.33"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("E50").Value = "  +3.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.27"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.52%  "
